$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update calculated value in row 2 (n changed) ---
$ws.Range("J2").Value = 65

# --- Update calculated values in row 3 ---
$ws.Range("B3").Value = 0.010488
$ws.Range("C3").Value = 0.004377
$ws.Range("D3").Value = 3.286204
$ws.Range("E3").Value = 0.191672
$ws.Range("F3").Value = 0.9344
$ws.Range("G3").Value = 0.8698
$ws.Range("H3").Value = 0.9993891
$ws.Range("I3").Value = 1.944467
$ws.Range("J3").Value = 56
$ws.Range("K3").Value = 8.065
$ws.Range("L3").Value = 17.748
$ws.Range("M3").Value = 75
$ws.Range("O3").Value = 95

# --- Update calculated values in row 4 ---
$ws.Range("B4").Value = 0.01451
$ws.Range("C4").Value = 0.004956
$ws.Range("D4").Value = 2.771572
$ws.Range("E4").Value = 0.207969
$ws.Range("F4").Value = 0.1844
$ws.Range("G4").Value = 0.7164
$ws.Range("H4").Value = 0.9991597
$ws.Range("I4").Value = 1.002258
$ws.Range("J4").Value = 124

# --- Update the active selection to match the saved view state ---
$ws.Range("G9").Select() | Out-Null

# --- Set page orientation to portrait (adds pageSetup element) ---
$ws.PageSetup.Orientation = 1 | Out-Null
